# 7.8 History Card & Advanced Story
# Update Fu's history-card dialogue lines on Sheet1 with the revised,
# more specific wording from the 7.8 script pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = "I used to live in JiuJiang county at the foot of the mountain."
$ws.Range("B5").Value = "Alas, my family was impoverished, and I possessed no particular skills——only a talent for poetry and calligraphy."
$ws.Range("B6").Value = "A month ago, I happened to meet Ming in downtown. We got along very well, and he invited me to reside at the manor."
$ws.Range("B8").Value = "When was the last time you saw the Lord?"
$ws.Range("B16").Value = "After around 3 PM, when most had likely finished their lunch, I went to the canteen."
$ws.Range("B18").Value = "Around 7.45 PM, Ming suddenly knocked on my door, asking if I knew where the Lord was."

# Restore the scrolled viewport / active selection recorded for this sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B23").Select()
